$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 219: add the missing "Category" cell (A219) ---------------------
# Re-use the formatting of A218 (same category "LIVE, TRAIN, RAIL") and
# just fill in the value.
$ws.Range("A218").Copy($ws.Range("A219"))
$ws.Cells.Item(219, 1).Value = "LIVE, TRAIN, RAIL"

# --- Row 220: Port Miami webcam (Florida, USA) ----------------------------
$ws.Range("A218:F218").Copy($ws.Range("A220:F220"))
$ws.Cells.Item(220, 1).Value = "LIVE, CITY, PORT, CRUSE, TRAFFIC"
$ws.Cells.Item(220, 2).Value = "25.775024065903935, -80.17222971340303"
$ws.Cells.Item(220, 3).Value = "LIVE - Port Miami Webcam with VHF Marine Radio Feed from PTZtv"
$ws.Cells.Item(220, 4).Value = "FL"
$ws.Cells.Item(220, 5).Value = "USA"
$ws.Cells.Item(220, 6).Value = "DxZziUUr6CY"

# --- Row 221: Miami Airport webcam (Florida, USA) -------------------------
$ws.Range("A218:F218").Copy($ws.Range("A221:F221"))
$ws.Cells.Item(221, 1).Value = "LIVE, AIRPORT"
$ws.Cells.Item(221, 2).Value = "25.792879308047315, -80.28243803922248"
$ws.Cells.Item(221, 3).Value = "LIVE Miami Airport -- Watch runway 9/27 with tower radio traffic!"
$ws.Cells.Item(221, 4).Value = "FL"
$ws.Cells.Item(221, 5).Value = "USA"
$ws.Cells.Item(221, 6).Value = "rDm2dFjRz3Q"

# --- Row 222: Taiwan Taoyuan International Airport ------------------------
$ws.Range("A217:F217").Copy($ws.Range("A222:F222"))
$ws.Cells.Item(222, 1).Value = "LIVE, AIRPORT"
$ws.Cells.Item(222, 2).Value = "25.065292854810945, 121.22902609240775"
$ws.Cells.Item(222, 3).Value = "Taiwan Taoyuan International Airport (TPE/RCTP) Live Camera 24/7"
$ws.Cells.Item(222, 4).Value = "Taoyuan"
$ws.Cells.Item(222, 5).Value = "Taiwan"
$ws.Cells.Item(222, 6).Value = "91PfFoqvuUk"

# --- Row 223: Hutoushan Environmental Park Live Cam ------------------------
$ws.Range("A217:F217").Copy($ws.Range("A223:F223"))
$ws.Cells.Item(223, 1).Value = "LIVE, CITY, LANDSCAPE"
$ws.Cells.Item(223, 2).Value = "25.0121745951993, 121.32839932553509"
$ws.Cells.Item(223, 3).Value = "Hutoushan Environmental Park Live Cam 桃園虎頭山環保公園即時影像"
$ws.Cells.Item(223, 4).Value = "Taoyuan"
$ws.Cells.Item(223, 5).Value = "Taiwan"
$ws.Cells.Item(223, 6).Value = "tu_gsIkNt-w"

# --- Update the view: active cell moves past the new last row -------------
$ws.Range("A224").Select()
